$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12, column A currently holds the phone number as text ("71277628").
# Convert it to a genuine numeric value, matching the target diff.
$ws.Cells.Item(12, 1).Value = 71277628

# Append a new redemption row (row 13) for the same phone number / points,
# a few seconds later. The phone number stays textual, like the rest of
# the sheet's "phone" column originally was.
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "71277628"
$ws.Cells.Item(13, 1).Style = "Normal"

$ws.Cells.Item(13, 2).Value = 76
$ws.Cells.Item(13, 3).Value = "2025-08-18T16:54:50"
